# Commit: "added icons to keys and minor formatting updates"
# Relevant spreadsheet-level changes:
#   1. Rename the "Mobotrex" sheet to "MoboTrex" (capitalization fix).
#   2. Make "MoboTrex" the active/selected sheet instead of "Western Systems".

$wb = $excel.ActiveWorkbook

# Fix the capitalization of the "Mobotrex" sheet name -> "MoboTrex"
$moboSheet = $wb.Worksheets.Item("Mobotrex")
$moboSheet.Name = "MoboTrex"

# Make the MoboTrex sheet the active tab (moves tabSelected from
# "Western Systems" to "MoboTrex" and updates workbookView activeTab).
$moboSheet.Activate()
